# Scheduled runner update: refresh market-price columns on the Leve profit
# sheets (ALC / BSM / LTW) — recompute one stale row on ALC and drop the
# cached price/profit figures (H:N) for the rows whose source data is no
# longer valid on BSM and LTW.

$wb = $excel.ActiveWorkbook

# --- ALC: row 51 gets recomputed averaged price/profit figures -------------
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3666.6667
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 3666.6667
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 3666.6667
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -4634.6667

# --- BSM: clear stale currentAveragePrice.../LevePrice.../LeveProfit... ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117:L120").ClearContents()
$ws.Range("H122:L135").ClearContents()
$ws.Range("M134").ClearContents()
$ws.Range("H137:L141").ClearContents()
$ws.Range("N137").ClearContents()

# --- LTW: same cleanup ------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124:L125").ClearContents()
$ws.Range("H127:L141").ClearContents()
$ws.Range("N127").ClearContents()
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()
